$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.090.25"
$ws.Range("E2").Value = "  +0.84%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.123.55"
$ws.Range("E3").Value = "  +1.45%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "577.63"
$ws.Range("E5").Value = "  -0.38%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "173.75"
$ws.Range("E6").Value = "  +3.14%  "
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.118.96"
$ws.Range("E8").Value = "  +1.39%  "
$ws.Range("E9").Value = "  -0.21%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.45"
$ws.Range("E10").Value = "  -3.93%  "
$ws.Range("E11").Value = "  +0.74%  "
$ws.Range("E12").Value = "  -1.43%  "
$ws.Range("E13").Value = "  -0.66%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "37.33"
$ws.Range("E14").Value = "  +1.04%  "
$ws.Range("E15").Value = "  -0.82%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.641.84"
$ws.Range("E16").Value = "  +2.02%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "67.113.91"
$ws.Range("E17").Value = "  +1.00%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.13"
$ws.Range("E18").Value = "  -1.72%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.126.89"
$ws.Range("E19").Value = "  +2.08%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "16.10"
$ws.Range("E20").Value = "  -3.12%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "483.45"
$ws.Range("E21").Value = "  +3.34%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.713"
$ws.Range("E22").Value = "  -0.20%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.68"
$ws.Range("E23").Value = "  +2.54%  "
$ws.Range("B24").Value = "Litecoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "83.97"
$ws.Range("E24").Value = "  +0.66%  "
$ws.Range("B25").Value = "InternetComputer(DFINITY)"
$ws.Range("C25").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "13.30"
$ws.Range("E25").Value = "  +3.03%  "
$ws.Range("E26").Value = "  +0.87%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.02"
$ws.Range("E27").Value = "  -0.86%  "
$ws.Range("E28").Value = "  -0.04%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.01"
$ws.Range("E29").Value = "  -2.25%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.40"
$ws.Range("E30").Value = "  -1.84%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.67"
$ws.Range("E31").Value = "  -0.03%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "28.71"
$ws.Range("E32").Value = "  +1.05%  "
$ws.Range("E33").Value = "  -2.76%  "
$ws.Range("E34").Value = "  -3.02%  "
$ws.Range("E35").Value = "  +0.22%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.87"
$ws.Range("E36").Value = "  -0.36%  "
$ws.Range("E37").Value = "  -1.61%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "47.74"
$ws.Range("E38").Value = "  -0.44%  "
$ws.Range("E39").Value = "  +1.88%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "50.11"
$ws.Range("E40").Value = "  +0.67%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.311"
$ws.Range("E41").Value = "  -2.85%  "
$ws.Range("E42").Value = "  +0.89%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.62"
$ws.Range("E43").Value = "  -1.09%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.64"
$ws.Range("E44").Value = "  -8.95%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.832.91"
$ws.Range("E45").Value = "  +2.83%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0357"
$ws.Range("E46").Value = "  -1.42%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "380.96"
$ws.Range("E47").Value = "  -0.82%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "136.06"
$ws.Range("E48").Value = "  +1.68%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "24.86"
$ws.Range("E50").Value = "  +0.85%  "
$ws.Range("E51").Value = "  -1.20%  "
